# "make check correct word"
#
# Slide 11 ("16. Бор") has two duplicate "Задача" textboxes sitting on top
# of each other (TextBox 2 / id=3 and TextBox 10 / id=11), each linking to
# a different LeetCode problem. The fix removes the duplicate (TextBox 10)
# and relabels the remaining textbox "Бор", carrying over the hyperlink
# that used to live on the run inside the deleted textbox.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(11)

$keep = $s.Shapes.Item("TextBox 2")
$dupe = $s.Shapes.Item("TextBox 10")

# Grab the hyperlink address from the run inside the duplicate box before
# it is removed, so it can be transplanted onto the surviving run.
$dupeRange = $dupe.TextFrame.TextRange
$dupeRun = $dupeRange.Characters(1, $dupeRange.Text.Length)
$targetAddress = $dupeRun.ActionSettings(1).Hyperlink.Address

$dupe.Delete()

$keep.TextFrame.TextRange.Text = "Бор"
$keepRange = $keep.TextFrame.TextRange
$keepRun = $keepRange.Characters(1, $keepRange.Text.Length)
$keepRun.ActionSettings(1).Hyperlink.Address = $targetAddress
